$d = $word.ActiveDocument

# Locate the italic "heure" run that sits alone in its own Heading2
# paragraph (the short "title-like" heading right below the "<" link
# paragraph). There is a second, non-italic "heure" Heading2 paragraph
# further down in the document which must be left untouched, so we
# search specifically for italic-formatted "heure" text.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Font.Italic = $true
$rng.Find.Text = "heure"
$found = $rng.Find.Execute()

if ($found) {
    # $rng now covers just the matched run ("heure"); expand it to the
    # whole enclosing paragraph (including the paragraph mark) so the
    # entire paragraph is removed, not just its text.
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            $target = $p
            break
        }
    }

    if ($target -ne $null) {
        $target.Range.Delete()
    }
}
